# Apply ticker-symbol updates to the "fidi" watchlist sheet (B2:F20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B
$ws.Range("B2").Value = "NSE:AARTIIND"
$ws.Range("B3").Value = "NSE:AMRUTANJAN"
$ws.Range("B4").Value = "NSE:APLAPOLLO"
$ws.Range("B5").Value = "NSE:BHAGERIA"
$ws.Range("B6").Value = "NSE:BHARATRAS"
$ws.Range("B7").Value = "NSE:CANFINHOME"
$ws.Range("B8").Value = "NSE:CIEINDIA"
$ws.Range("B9").Value = "NSE:DATAMATICS"
$ws.Range("B10").Value = "NSE:GENESYS"
$ws.Range("B11").Value = "NSE:GOODLUCK"
$ws.Range("B12").Value = "NSE:GUJGASLTD"
$ws.Range("B13").Value = "NSE:HEADSUP"
$ws.Range("B14").Value = "NSE:HEXATRADEX"
$ws.Range("B15").Value = "NSE:IGL"
$ws.Range("B16").Value = "NSE:IPCALAB"
$ws.Range("B17").Value = "NSE:LINC"
$ws.Range("B18").Value = "NSE:MOLDTECH"
$ws.Range("B19").Value = "NSE:NAHARCAP"
$ws.Range("B20").Value = "NSE:ROSSELLIND"

# Column C
$ws.Range("C2").Value = "NSE:APLLTD"
$ws.Range("C3").Value = "NSE:ATL"
$ws.Range("C4").Value = "NSE:BLISSGVS"
$ws.Range("C5").Value = "NSE:CALSOFT"
$ws.Range("C6").Value = "NSE:GKWLIMITED"
$ws.Range("C7").Value = "NSE:GREENLAM"
$ws.Range("C8").Value = "NSE:GRINFRA"
$ws.Range("C9").Value = "NSE:HINDCON"
$ws.Range("C10").Value = "NSE:HMVL"
$ws.Range("C11").Value = "NSE:KIRLPNU"
$ws.Range("C12").Value = "NSE:NILAINFRA"
$ws.Range("C13").Value = "NSE:ORIENTBELL"
$ws.Range("C14").Value = "NSE:REDTAPE"
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()

# Column D
$ws.Range("D2").Value = "NSE:IPCALAB"
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()

# Column E
$ws.Range("E3").Value = "NSE:HDFCBANK"
$ws.Range("E4").Value = "NSE:KOTAKBANK"

# Column F
$ws.Range("F2").Value = "NSE:BAJAJ-AUTO"
$ws.Range("F3").Value = "NSE:IGL"
$ws.Range("F4").Value = "NSE:IPCALAB"

# Rows 21-26 no longer used; remove them so the sheet dimension becomes A1:F20
$ws.Range("A21:F26").EntireRow.Delete()
